# Entrega Introducción Bety Díaz actividad 4 - apply commit edits
#
# Real content changes identified from the OOXML diff (everything else in
# the diff is just <w:lastRenderedPageBreak/> churn caused by Word
# re-paginating after the text insertion below - no actual text changed
# there, so it is left for Word to recompute naturally):
#
#   1. "...entre otros resultados encontró que debido..." ->
#      "...entre otros resultados encontró que, debido..." (adds a comma;
#      the grammar-check proofErr wrapper around the old "que" is dropped
#      along with it).
#   2. The quantitative-strategies paragraph gets a new block of sentences
#      appended after "...fase cuantitativa son: ".

$d = $word.ActiveDocument

# 1) "que" -> "que,"
$null = $d.Content.Find.Execute(
    "encontró que debido",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "encontró que, debido",
    2)

# 2) Append the new instruments/strategies sentences to the paragraph that
#    currently ends right after "...fase cuantitativa son: ".
$rng = $d.Content
$null = $rng.Find.Execute(
    "Las estrategias de recolección y análisis de datos de la fase cuantitativa son: ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "",
    0)
$rng.Collapse(0)
$rng.InsertAfter("Análisis de Regresión Multivariable. Las estrategias de recolección de datos son: entrevistas estructuradas y test estandarizados. Los instrumentos son: General – HVD: Indicadores de Bienestar y Salud Mental en el Mundo del Trabajo de la Universidad de La Sabana; Cuestionario de Estrés Laboral (CESQT); Escala de Ansiedad de Hamilton (HAM-A); Cuestionario de Depresión en el Trabajo (WDQ - Workplace Depression Questionnaire).")
